# Commit: Client 에서 Cookie Enhance 호출
# The two rows describing BefAccSoulStone / AftAccSoulStone (the old
# C91/C92 values under CookieEnhanceStarReq) are removed from the
# "Packet" sheet. Every row below shifts up by two; Excel automatically
# prunes the now-unused shared strings (BefAccSoulStone, AftAccSoulStone)
# from the shared string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")
$ws.Activate()

# Delete the two rows (Excel rows 91 and 92) that hold the
# BefAccSoulStone / AftAccSoulStone packet fields.
$ws.Rows("91:92").Delete()

# Restore the selection/view to where the author ended up after the edit.
$ws.Range("C90").Select()
